$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was refreshed
# for every data row (rows 2-264) from 45192 (2023-09-23) to 45202 (2023-10-03).
$ws.Range("C2:C264").Value = 45202
